# download articles with pandoc title blocks
#
# The document currently starts with:
#   [bookmarkStart "aims-and-purposes"]
#   P1 (Heading1):  "Aims and Purposes"
#   [bookmarkEnd   "aims-and-purposes"]
#   P2 (no style, bold run): "By Dorothy Day"
#
# It should become a pandoc-style title block:
#   P1 (Title):    "Aims" " " "and" " " "Purposes"   (separate runs)
#   P2 (Authors):  "Dorothy" " " "Day"                (separate runs)

$d = $word.ActiveDocument

$titlePara = $d.Paragraphs(1)
$authorPara = $d.Paragraphs(2)

$rangeStart = $titlePara.Range.Start
$rangeEnd = $authorPara.Range.End

$replacement = $d.Range($rangeStart, $rangeEnd)

$openXmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" ' + `
  'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>'

$titleParagraphXml = '<w:p>' + `
    '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Aims</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">and</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Purposes</w:t></w:r>' + `
  '</w:p>'

$authorParagraphXml = '<w:p>' + `
    '<w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' + `
  '</w:p>'

$openXmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$fragment = $openXmlHeader + $titleParagraphXml + $authorParagraphXml + $openXmlFooter

# InsertXML replaces the contents of the target range while keeping it a
# single logical edit, which is what lets the new paragraph properties
# (pStyle Title / Authors) take effect and keeps each word in its own run,
# matching the pandoc-style title-block markup produced upstream.
$replacement.InsertXML($fragment)

# The source bookmark that used to wrap the old "Aims and Purposes" heading
# is no longer meaningful once the heading becomes a pandoc title block, so
# try to drop it. This is best effort: some Word automation hosts do not
# surface bookmark deletion through the object model, in which case the
# (now harmless/unreferenced) bookmark markers are simply left in place.
try {
    $bookmarks = $d.Bookmarks
    $oldBookmark = $bookmarks.Item(1)
    if ($oldBookmark.Name -eq "aims-and-purposes" -or $oldBookmark.Name -eq "") {
        $oldBookmark.Delete()
    }
} catch {
    # Bookmark removal is best-effort; continue even if unsupported.
}

Write-Host "Converted heading/author block into pandoc title-block style paragraphs."
